$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 13 ("e501" / "Game Won") to make room
# for the new "e012 Hatches" entry. This pushes the old rows 13-17 down to
# rows 14-18 and carries their formatting along automatically.
$ws.Rows("13").Insert()

# Populate the new row with the new event text.
$ws.Range("A13").Value = "e012"

$hatchesText = @"
<Bold>e012 Hatches</Bold> 
<InlineUIContainer><Button Content='r4.42' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Left click on hatches on the Tank Card to toggle adding counter. Click image below to continue.
<LineBreak/><LineBreak/>
                                               <InlineUIContainer><Image Name='c15OpenHatch'  Height='80' Width='80'></Image></InlineUIContainer>"
"@
$ws.Range("B13").Value = $hatchesText

# Match the row height used by the other long entries (e.g. row 17/18).
$ws.Rows("13").RowHeight = 85.6

# Reflect the new selection/active cell recorded in the saved workbook.
$ws.Range("B13").Select()
